# Apply the "week 10 / Wed 2018.11.7" update to the project plan sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- 1) Mark the six previous-week tasks (rows 133-138) as completed ---
$ws.Range("C133").Value = "已完成"
$ws.Range("C134").Value = "已完成"
$ws.Range("C135").Value = "已完成"
$ws.Range("C136").Value = "已完成"
$ws.Range("C137").Value = "已完成"
$ws.Range("C138").Value = "已完成"

# Row 135's task description changes from the e-mail verification code task
# to the login picture-verification-code task.
$ws.Range("B135").Value = "android登陆图片验证码功能实现"

# The closing summary for that week now gets a fuller write-up.
$ws.Range("A139").Value = "总结：前端成功集成了环信SDK，重心移向后台。"

# --- 2) Append a brand-new week section: rows 141-150 ---
# Merge the new banner/summary blocks *before* copying any formatting onto
# them (merging after formatting makes Excel redistribute the outer-box
# border across the merged cells' individual edges, which does not match
# how the rest of this sheet stores merged-row styles: every cell in a
# merged run shares one plain style index).
$ws.Range("A141:D141").Merge()
$ws.Range("A149:D150").Merge()

# Now stamp the row styles by copying *formats only* from the matching
# rows of the previous week's block (131-140), so fonts/borders/number
# formats line up with the sheet's existing section pattern:
#   banner row -> header row -> 6 data rows -> summary block
$ws.Range("A131:D131").Copy()
$ws.Range("A141").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = 0

$ws.Range("A132:D138").Copy()
$ws.Range("A142").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = 0

$ws.Range("A139:D140").Copy()
$ws.Range("A149").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = 0

# Date banner for the new week.
$ws.Range("A141").Value = "日期：2018.11.7 第十周周三"

# Column headers.
$ws.Range("A142").Value = "组员"
$ws.Range("B142").Value = "计划内容"
$ws.Range("C142").Value = "完成情况"
$ws.Range("D142").Value = "备注"

# Six members' planned tasks for the new week.
$ws.Range("A143").Value = "练富珊"
$ws.Range("B143").Value = "利用环信实现聊天功能"

$ws.Range("A144").Value = "黄成志"
$ws.Range("B144").Value = "后台数据库建表"

$ws.Range("A145").Value = "黄皓燊"
$ws.Range("B145").Value = "android忘记密码手机验证码功能实现"

$ws.Range("A146").Value = "郑嘉蔚"
$ws.Range("B146").Value = "制作静态后台管理网页[开启群]"

$ws.Range("A147").Value = "陈碧容"
$ws.Range("B147").Value = "制作静态后台管理网页[禁用用户]"

$ws.Range("A148").Value = "辛伟达"
$ws.Range("B148").Value = "制作静态后台管理网页[禁用群]"

# Trailing (still-empty) summary block for the new week.
$ws.Range("A149").Value = "总结："

# --- 3) Update the sheet view bookkeeping to match ---
$ws.Range("A124").Select()
$excel.ActiveWindow.ScrollRow = 124
$ws.Range("E141").Select()
